$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 23 (Entergy Transmission - SOC, Louisiana) ---
# The restoration date/duration, previously "ongoing", are now reported with
# placeholder dot-filled values.
$ws.Range("D23").Value = ".        ."
$ws.Range("E23").Value = ". Hours,  . Minutes"

# --- Insert a new disturbance record as row 83 (Modesto Irrigation District) ---
# This pushes the existing footnote row (old row 83) down to row 84.
$ws.Rows(83).Insert()

# Copy formatting from the row above (an existing data row) onto the new row
# so the new row gets the same style (borders, alignment, font) as the other
# data rows instead of Excel's blank-insert default.
$ws.Range("A82:K82").Copy()
$ws.Range("A83:K83").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A83").Value = 2016
$ws.Range("B83").Value = 11
$ws.Range("C83").Value = "11/09/2016 11:59 AM"
$ws.Range("D83").Value = "11/09/2016  6:15 PM"
$ws.Range("E83").Value = "6 Hours, 16 Minutes"
$ws.Range("F83").Value = "Modesto Irrigation District"
$ws.Range("G83").Value = "WECC"
$ws.Range("H83").Value = "California: Stanislaus County, San Joaquin County, Alameda County, Tuolumne County;"
$ws.Range("I83").Value = "Cyber event that could potentially impact electric power system adequacy or reliability-Cyber Attack"
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0

# Match the row height Excel computed for this wrapped text in the source file.
$ws.Rows(83).RowHeight = 51.75

# --- Restore the active-cell selection left by the editor ---
[void]$ws.Range("E29").Select()
